$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Cells.Item(20, 1).Value() = 131199758
$ws.Cells.Item(20, 2).Value() = 5197
$ws.Cells.Item(20, 4).Value() = "LC"
$ws.Cells.Item(20, 5).Value() = 105930
$ws.Cells.Item(20, 6).Value() = "Vågbandad barkbock"
$ws.Cells.Item(20, 7).Value() = "Semanotus undatus"
$ws.Cells.Item(20, 8).Value() = "(Linnaeus, 1758)"
# $ws.Cells.Item(20, 9) left blank (empty string in source)
$ws.Cells.Item(20, 13).Value() = "äldre gnagspår"
$ws.Cells.Item(20, 16).Value() = "Anderstorp, Anderstorp, Nrk"
$ws.Cells.Item(20, 17).Value() = 467630
$ws.Cells.Item(20, 18).Value() = 6544255
$ws.Cells.Item(20, 19).Value() = 20
$ws.Cells.Item(20, 20).Value() = "Örebro"
$ws.Cells.Item(20, 21).Value() = "Degerfors"
$ws.Cells.Item(20, 22).Value() = "Närke"
$ws.Cells.Item(20, 23).Value() = "Nysund"
$ws.Cells.Item(20, 25).NumberFormat = "@"
$ws.Cells.Item(20, 25).Value() = "2026-02-17"
$ws.Cells.Item(20, 26).Value() = "13:27"
$ws.Cells.Item(20, 27).NumberFormat = "@"
$ws.Cells.Item(20, 27).Value() = "2026-02-17"
$ws.Cells.Item(20, 28).Value() = "13:27"
$ws.Cells.Item(20, 29).Value() = "Äldre gnagspår och kläckhål på stående död gran"
$ws.Cells.Item(20, 30).Value() = $false
$ws.Cells.Item(20, 31).Value() = $false
$ws.Cells.Item(20, 33).Value() = $false
# $ws.Cells.Item(20, 46) left blank (empty string in source)
$ws.Cells.Item(20, 49).Value() = "Therese Steiner"
$ws.Cells.Item(20, 50).Value() = "Therese Steiner"
# $ws.Cells.Item(20, 51) left blank (empty string in source)

# Row 21
$ws.Cells.Item(21, 1).Value() = 131199884
$ws.Cells.Item(21, 2).Value() = 4773
$ws.Cells.Item(21, 4).Value() = "LC"
$ws.Cells.Item(21, 5).Value() = 100299
$ws.Cells.Item(21, 6).Value() = "Thomsons trägnagare"
$ws.Cells.Item(21, 7).Value() = "Cacotemnus thomsoni"
$ws.Cells.Item(21, 8).Value() = "(Kraatz, 1881)"
# $ws.Cells.Item(21, 9) left blank (empty string in source)
$ws.Cells.Item(21, 16).Value() = "Anderstorp, Anderstorp, Nrk"
$ws.Cells.Item(21, 17).Value() = 467729
$ws.Cells.Item(21, 18).Value() = 6544345
$ws.Cells.Item(21, 19).Value() = 20
$ws.Cells.Item(21, 20).Value() = "Örebro"
$ws.Cells.Item(21, 21).Value() = "Degerfors"
$ws.Cells.Item(21, 22).Value() = "Närke"
$ws.Cells.Item(21, 23).Value() = "Nysund"
$ws.Cells.Item(21, 25).NumberFormat = "@"
$ws.Cells.Item(21, 25).Value() = "2026-02-17"
$ws.Cells.Item(21, 26).Value() = "13:37"
$ws.Cells.Item(21, 27).NumberFormat = "@"
$ws.Cells.Item(21, 27).Value() = "2026-02-17"
$ws.Cells.Item(21, 28).Value() = "13:37"
$ws.Cells.Item(21, 29).Value() = "Gnagspår och kläckhål i stående död gran"
$ws.Cells.Item(21, 30).Value() = $false
$ws.Cells.Item(21, 31).Value() = $false
$ws.Cells.Item(21, 33).Value() = $false
# $ws.Cells.Item(21, 46) left blank (empty string in source)
$ws.Cells.Item(21, 49).Value() = "Therese Steiner"
$ws.Cells.Item(21, 50).Value() = "Therese Steiner"
# $ws.Cells.Item(21, 51) left blank (empty string in source)

